$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (serial 45202 = 2023-10-03)
# for every data row (2..92); bump it by one day to serial 45203 (2023-10-04).
for ($r = 2; $r -le 92; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}
